$p = $ppt.ActivePresentation

# Slide 4 ("Update" - Implemented/Back Testing in progress/Planned) currently has
# no speaker notes. Add a notes page with the back-testing definition, matching
# the bold lead-in "Back testing" followed by the explanatory sentence.
$s4 = $p.Slides.Item(4)
$notesShape = $s4.NotesPage.Shapes.AddPlaceholder(2)
$notesShape.TextFrame.TextRange.Text = "Back testing is the process of evaluating a trading strategy or investment model by applying it to historical data to see how it would have performed in the past."
